$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a number/numeric string: force text format
# so Excel stores the literal text (preserving formatting like trailing zeros,
# leading zeros, exact precision) instead of silently converting to a float.
function Set-TextValue($rangeAddr, $value) {
    $rng = $ws.Range($rangeAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" '249.69'
Set-TextValue "G2" '22'
Set-TextValue "D3" '21.64'
Set-TextValue "G3" '22'
Set-TextValue "D4" '5.583'
Set-TextValue "G4" '22'
Set-TextValue "D5" '0.05672'
Set-TextValue "G5" '22'
Set-TextValue "D6" '3.377'
Set-TextValue "G6" '22'
Set-TextValue "D7" '6.432'
Set-TextValue "G7" '22'
Set-TextValue "D8" '0.8060'
Set-TextValue "G8" '22'
Set-TextValue "D9" '1.040'
Set-TextValue "G9" '22'
Set-TextValue "D10" '0.1426'
Set-TextValue "G10" '22'
Set-TextValue "D11" '0.07250'
Set-TextValue "G11" '22'
Set-TextValue "D12" '0.03124'
Set-TextValue "G12" '22'
Set-TextValue "D13" '0.02922'
Set-TextValue "G13" '22'
Set-TextValue "D14" '0.09268'
Set-TextValue "G14" '22'
Set-TextValue "D15" '0.001669'
Set-TextValue "G15" '22'
Set-TextValue "G16" '22'
Set-TextValue "D17" '0.04725'
Set-TextValue "G17" '22'
Set-TextValue "D18" '0.0005814'
Set-TextValue "G18" '22'
Set-TextValue "D19" '0.006460'
Set-TextValue "G19" '22'
Set-TextValue "D20" '0.005056'
Set-TextValue "G20" '22'
Set-TextValue "D21" '0.001051'
Set-TextValue "G21" '22'
Set-TextValue "G22" '22'
Set-TextValue "D23" '3.981'
Set-TextValue "G23" '22'
Set-TextValue "D24" '2.114'
Set-TextValue "G24" '22'
Set-TextValue "G25" '22'
Set-TextValue "G26" '22'
Set-TextValue "D27" '0.0003102'
Set-TextValue "G27" '22'
Set-TextValue "G28" '22'
Set-TextValue "G29" '22'
Set-TextValue "G30" '22'
Set-TextValue "G31" '22'
Set-TextValue "G32" '22'
Set-TextValue "G33" '22'
Set-TextValue "G34" '22'
Set-TextValue "G35" '22'
Set-TextValue "G36" '22'
Set-TextValue "G37" '22'
Set-TextValue "G38" '22'
Set-TextValue "G39" '22'
Set-TextValue "D40" '0.04125'
Set-TextValue "G40" '22'
Set-TextValue "D41" '0.006904'
Set-TextValue "G41" '22'
$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
Set-TextValue "D42" '0.1043'
$ws.Range("E42").Value = '41BKEXTokenBKK'
Set-TextValue "G42" '22'
$ws.Range("B43").Value = 'CEJI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
Set-TextValue "D43" '0.002972'
$ws.Range("E43").Value = '42CEJICEJI'
Set-TextValue "G43" '22'
Set-TextValue "D44" '0.008533'
Set-TextValue "G44" '22'
Set-TextValue "G45" '22'
Set-TextValue "G46" '22'
Set-TextValue "D47" '0.7857'
Set-TextValue "G47" '22'
Set-TextValue "D48" '0.01666'
Set-TextValue "G48" '22'
Set-TextValue "D49" '0.00002101'
Set-TextValue "G49" '22'
Set-TextValue "G50" '22'
Set-TextValue "G51" '22'
